{"js": "const replacements = [\n  [\"2024-01-24 Wednesday\", \"2024-01-25 Thursday\"],\n  [\"33\u00d720=\", \"55\u00d720=\"],\n  [\"24\u00d744=\", \"64\u00d780=\"],\n  [\"65\u00d762=\", \"76\u00d792=\"],\n  [\"40\u00d753=\", \"83\u00d737=\"],\n  [\"92\u00d798=\", \"28\u00d749=\"],\n  [\"53\u00d764=\", \"44\u00d718=\"],\n  [\"38\u00d793=\", \"63\u00d792=\"],\n  [\"84\u00d779=\", \"18\u00d757=\"],\n  [\"99\u00d739=\", \"15\u00d771=\"],\n  [\"75\u00d764=\", \"37\u00d749=\"],\n  [\"75\u00d796=\", \"25\u00d726=\"],\n  [\"59\u00d774=\", \"64\u00d731=\"],\n  [\"54\u00d799=\", \"89\u00d759=\"],\n  [\"26\u00d735=\", \"76\u00d758=\"],\n  [\"54\u00d741=\", \"22\u00d731=\"],\n  [\"84\u00d724=\", \"26\u00d756=\"],\n  [\"30\u00d773=\", \"46\u00d735=\"],\n  [\"49\u00d793=\", \"56\u00d721=\"],\n  [\"62\u00d754=\", \"24\u00d754=\"],\n  [\"40\u00d787=\", \"39\u00d743=\"],\n  [\"36\u00d765=\", \"71\u00d724=\"],\n  [\"43\u00d757=\", \"48\u00d788=\"],\n  [\"18\u00d755=\", \"45\u00d751=\"],\n  [\"18\u00d793=\", \"59\u00d785=\"],\n  [\"28\u00d720=\", \"51\u00d734=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @('2024-01-24 Wednesday', '2024-01-25 Thursday'),\n    @('33\u00d720=', '55\u00d720='),\n    @('24\u00d744=', '64\u00d780='),\n    @('65\u00d762=', '76\u00d792='),\n    @('40\u00d753=', '83\u00d737='),\n    @('92\u00d798=', '28\u00d749='),\n    @('53\u00d764=', '44\u00d718='),\n    @('38\u00d793=', '63\u00d792='),\n    @('84\u00d779=', '18\u00d757='),\n    @('99\u00d739=', '15\u00d771='),\n    @('75\u00d764=', '37\u00d749='),\n    @('75\u00d796=', '25\u00d726='),\n    @('59\u00d774=', '64\u00d731='),\n    @('54\u00d799=', '89\u00d759='),\n    @('26\u00d735=', '76\u00d758='),\n    @('54\u00d741=', '22\u00d731='),\n    @('84\u00d724=', '26\u00d756='),\n    @('30\u00d773=', '46\u00d735='),\n    @('49\u00d793=', '56\u00d721='),\n    @('62\u00d754=', '24\u00d754='),\n    @('40\u00d787=', '39\u00d743='),\n    @('36\u00d765=', '71\u00d724='),\n    @('43\u00d757=', '48\u00d788='),\n    @('18\u00d755=', '45\u00d751='),\n    @('18\u00d793=', '59\u00d785='),\n    @('28\u00d720=', '51\u00d734='),\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    [void]$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 0, $false, $find.Replacement.Text, 2)\n}\n"}
